# "error solve ifrs list" - correct the financial figures for 쿠쿠홀딩스
# (values had been scaled/merged incorrectly); also drop the stale
# yearly columns/rows that no longer have corresponding data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 5667
$ws.Range("E2").Value = 786
$ws.Range("F2").Value = 786
$ws.Range("G2").Value = 1180
$ws.Range("H2").Value = 907
$ws.Range("I2").Value = 907
$ws.Range("K2").Value = 5533
$ws.Range("L2").Value = 1124
$ws.Range("M2").Value = 4410
$ws.Range("N2").Value = 4410
$ws.Range("P2").Value = 49
$ws.Range("Q2").Value = 576
$ws.Range("R2").Value = -526
$ws.Range("S2").Value = -41
$ws.Range("T2").Value = 191
$ws.Range("U2").Value = 386
$ws.Range("V2").Value = 30
$ws.Range("W2").Value = 13.87
$ws.Range("X2").Value = 16.01
$ws.Range("Y2").Value = 22.49
$ws.Range("Z2").Value = 17.98
$ws.Range("AA2").Value = 25.49
$ws.Range("AB2").Value = 10433.73
$ws.Range("AC2").Value = 9255
$ws.Range("AD2").Value = 17.31
$ws.Range("AE2").Value = 53790
$ws.Range("AF2").Value = 2.98
$ws.Range("AG2").Value = 1500
$ws.Range("AH2").Value = 0.9399999999999999
$ws.Range("AI2").Value = 13.55
$ws.Range("AJ2").Value = 9803360
$ws.Range("J2").ClearContents()
$ws.Range("O2").ClearContents()

# Row 3
$ws.Range("D3").Value = 6675
$ws.Range("E3").Value = 916
$ws.Range("F3").Value = 916
$ws.Range("G3").Value = 979
$ws.Range("H3").Value = 746
$ws.Range("I3").Value = 746
$ws.Range("J3").Value = -1
$ws.Range("K3").Value = 6210
$ws.Range("L3").Value = 1148
$ws.Range("M3").Value = 5062
$ws.Range("N3").Value = 5053
$ws.Range("O3").Value = 9
$ws.Range("P3").Value = 49
$ws.Range("Q3").Value = 211
$ws.Range("R3").Value = -38
$ws.Range("S3").Value = -153
$ws.Range("T3").Value = 190
$ws.Range("U3").Value = 21
$ws.Range("V3").Value = 1
$ws.Range("W3").Value = 13.73
$ws.Range("X3").Value = 11.17
$ws.Range("Y3").Value = 15.77
$ws.Range("Z3").Value = 12.7
$ws.Range("AA3").Value = 22.68
$ws.Range("AB3").Value = 11718.06
$ws.Range("AC3").Value = 7611
$ws.Range("AD3").Value = 27.14
$ws.Range("AE3").Value = 61633
$ws.Range("AF3").Value = 3.35
$ws.Range("AG3").Value = 2100
$ws.Range("AH3").Value = 1.02
$ws.Range("AI3").Value = 23.07
$ws.Range("AJ3").Value = 9803360

# Row 4
$ws.Range("D4").Value = 4838
$ws.Range("E4").Value = 638
$ws.Range("F4").Value = 954
$ws.Range("G4").Value = 791
$ws.Range("H4").Value = 801
$ws.Range("I4").Value = 811
$ws.Range("J4").Value = -10
$ws.Range("K4").Value = 7148
$ws.Range("L4").Value = 1475
$ws.Range("M4").Value = 5673
$ws.Range("N4").Value = 5641
$ws.Range("O4").Value = 32
$ws.Range("P4").Value = 49
$ws.Range("Q4").Value = 187
$ws.Range("R4").Value = -69
$ws.Range("S4").Value = -140
$ws.Range("T4").Value = 180
$ws.Range("U4").Value = 7
$ws.Range("V4").Value = 1
$ws.Range("W4").Value = 13.2
$ws.Range("X4").Value = 16.55
$ws.Range("Y4").Value = 15.17
$ws.Range("Z4").Value = 11.99
$ws.Range("AA4").Value = 25.99
$ws.Range("AB4").Value = 13018.45
$ws.Range("AC4").Value = 8272
$ws.Range("AD4").Value = 14.04
$ws.Range("AE4").Value = 69180
$ws.Range("AF4").Value = 1.68
$ws.Range("AG4").Value = 3100
$ws.Range("AH4").Value = 2.67
$ws.Range("AI4").Value = 31.17
$ws.Range("AJ4").Value = 9803360

# Row 5
$ws.Range("D5").Value = 4502
$ws.Range("E5").Value = 647
$ws.Range("F5").Value = 647
$ws.Range("G5").Value = 749
$ws.Range("H5").Value = 4520
$ws.Range("I5").Value = 4470
$ws.Range("J5").Value = 50
$ws.Range("K5").Value = 4896
$ws.Range("L5").Value = 1318
$ws.Range("M5").Value = 3578
$ws.Range("N5").Value = 3578
$ws.Range("P5").Value = 27
$ws.Range("Q5").Value = -301
$ws.Range("R5").Value = 367
$ws.Range("S5").Value = -363
$ws.Range("T5").Value = 134
$ws.Range("U5").Value = -435
$ws.Range("V5").Value = 0
$ws.Range("W5").Value = 14.37
$ws.Range("X5").Value = 100.4
$ws.Range("Y5").Value = 96.95999999999999
$ws.Range("Z5").Value = 75.05
$ws.Range("AA5").Value = 36.84
$ws.Range("AB5").Value = 38836.22
$ws.Range("AC5").Value = 47436
$ws.Range("AD5").Value = 3.13
$ws.Range("AE5").Value = 80971
$ws.Range("AF5").Value = 1.83
$ws.Range("AG5").Value = 4100
$ws.Range("AH5").Value = 2.76
$ws.Range("AI5").Value = 4.05
$ws.Range("AJ5").Value = 5315894
$ws.Range("O5").ClearContents()

# Row 6
$ws.Range("D6").Value = 4933
$ws.Range("E6").Value = 679
$ws.Range("F6").Value = 679
$ws.Range("G6").Value = 1045
$ws.Range("H6").Value = 846
$ws.Range("I6").Value = 846
$ws.Range("K6").Value = 7470
$ws.Range("L6").Value = 1376
$ws.Range("M6").Value = 6095
$ws.Range("N6").Value = 6095
$ws.Range("P6").Value = 36
$ws.Range("Q6").Value = 731
$ws.Range("R6").Value = -502
$ws.Range("S6").Value = -184
$ws.Range("T6").Value = 35
$ws.Range("U6").Value = 696
$ws.Range("V6").Value = 0
$ws.Range("W6").Value = 13.76
$ws.Range("X6").Value = 17.14
$ws.Range("Y6").Value = 17.49
$ws.Range("Z6").Value = 13.68
$ws.Range("AA6").Value = 22.57
$ws.Range("AB6").Value = 36135.32
$ws.Range("AC6").Value = 13096
$ws.Range("AD6").Value = 9.050000000000001
$ws.Range("AE6").Value = 98052
$ws.Range("AF6").Value = 1.21
$ws.Range("AG6").Value = 3000
$ws.Range("AH6").Value = 2.53
$ws.Range("AI6").Value = 22.05
$ws.Range("AJ6").Value = 7112437

# Rows 7-9: clear all data columns (D:AJ), keep A/B/C
$ws.Range("D7:AJ9").ClearContents()
